$d = $word.ActiveDocument

# --- Change 1: title "МДНФ и МКНФ" -> "СДНФ и СКНФ" (kept italic, same run formatting) ---
$d.Content.Find.Execute("МДНФ и МКНФ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "СДНФ и СКНФ", 2)

# --- Change 2: fix spelling "ассисент" -> "ассистент" ---
$d.Content.Find.Execute("ассисент", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ассистент", 2)
